# This script applies updated betting-odds values to Sheet1 of the workbook,
# matching the changes described in the commit "Atualizando o arquivo XLSX".
# Most edits update existing numeric odds; rows 12 and 13 (Latvia Virsliga
# fixtures) had most of their odds columns populated for the first time
# (they were previously blank placeholder cells), while columns J/K (and,
# for row 13, P/Q) remain blank as in the original file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("L2").Value = 1.44
$ws.Range("M2").Value = 2.63
$ws.Range("N2").Value = 2.35
$ws.Range("O2").Value = 1.57
# Row 3
$ws.Range("J3").Value = 1.1
$ws.Range("K3").Value = 7
# Row 8
$ws.Range("G8").Value = 2.85
$ws.Range("H8").Value = 2.5
$ws.Range("I8").Value = 3.05
$ws.Range("J8").Value = 1.19
$ws.Range("K8").Value = 4.1
$ws.Range("L8").Value = 1.78
$ws.Range("M8").Value = 1.93
$ws.Range("N8").Value = 3.25
$ws.Range("O8").Value = 1.29
$ws.Range("P8").Value = 1.75
$ws.Range("Q8").Value = 1.98
$ws.Range("R8").Value = 2.42
$ws.Range("S8").Value = 1.5
$ws.Range("T8").Value = 5.6
$ws.Range("V8").Value = 11.75
$ws.Range("Y8").Value = 65
$ws.Range("Z8").Value = 4.1
$ws.Range("AB8").Value = 22
$ws.Range("AC8").Value = 200
$ws.Range("AD8").Value = 5.6
$ws.Range("AE8").Value = 13
$ws.Range("AF8").Value = 13
$ws.Range("AG8").Value = 45
$ws.Range("AH8").Value = 45
$ws.Range("AI8").Value = 80
# Row 9
$ws.Range("G9").Value = 1.66
$ws.Range("J9").Value = 1.04
$ws.Range("L9").Value = 1.3
$ws.Range("O9").Value = 1.72
# Row 10
$ws.Range("G10").Value = 1.37
$ws.Range("H10").Value = 4.3
$ws.Range("I10").Value = 7.7
$ws.Range("L10").Value = 1.21
$ws.Range("M10").Value = 3.55
$ws.Range("N10").Value = 1.65
$ws.Range("O10").Value = 2
$ws.Range("R10").Value = 1.87
$ws.Range("S10").Value = 1.75
$ws.Range("T10").Value = 6.9
$ws.Range("U10").Value = 6.5
$ws.Range("W10").Value = 8.75
$ws.Range("X10").Value = 11.25
$ws.Range("Y10").Value = 26
$ws.Range("Z10").Value = 12.5
$ws.Range("AA10").Value = 8.75
$ws.Range("AB10").Value = 19
$ws.Range("AC10").Value = 90
$ws.Range("AD10").Value = 21
$ws.Range("AE10").Value = 55
$ws.Range("AF10").Value = 24
$ws.Range("AI10").Value = 70
$ws.Range("AJ10").Value = 700
# Row 12
$ws.Range("G12").Value = 2.77
$ws.Range("H12").Value = 3.15
$ws.Range("I12").Value = 2.32
$ws.Range("L12").Value = 1.28
$ws.Range("M12").Value = 3.3
$ws.Range("N12").Value = 1.85
$ws.Range("O12").Value = 1.75
$ws.Range("P12").Value = 1.38
$ws.Range("Q12").Value = 2.45
$ws.Range("R12").Value = 1.72
$ws.Range("S12").Value = 2.01
$ws.Range("T12").Value = 7.8
$ws.Range("U12").Value = 12
$ws.Range("V12").Value = 8.5
$ws.Range("W12").Value = 26
$ws.Range("X12").Value = 18.5
$ws.Range("Y12").Value = 24
$ws.Range("Z12").Value = 9.25
$ws.Range("AA12").Value = 5.4
$ws.Range("AB12").Value = 10.75
$ws.Range("AC12").Value = 45
$ws.Range("AD12").Value = 6.9
$ws.Range("AE12").Value = 9.75
$ws.Range("AF12").Value = 7.7
$ws.Range("AG12").Value = 19
$ws.Range("AH12").Value = 15.5
$ws.Range("AI12").Value = 22
$ws.Range("AJ12").Value = 300
# Row 13
$ws.Range("G13").Value = 1.87
$ws.Range("H13").Value = 3.55
$ws.Range("I13").Value = 3.4
$ws.Range("L13").Value = 1.25
$ws.Range("M13").Value = 3.6
$ws.Range("N13").Value = 1.75
$ws.Range("O13").Value = 1.87
$ws.Range("R13").Value = 1.71
$ws.Range("S13").Value = 2.03
$ws.Range("T13").Value = 6.7
$ws.Range("U13").Value = 7.8
$ws.Range("V13").Value = 7.2
$ws.Range("W13").Value = 13
$ws.Range("X13").Value = 12
$ws.Range("Y13").Value = 19.5
$ws.Range("Z13").Value = 11.25
$ws.Range("AA13").Value = 6.1
$ws.Range("AB13").Value = 11.75
$ws.Range("AC13").Value = 45
$ws.Range("AD13").Value = 9.5
$ws.Range("AE13").Value = 15.5
$ws.Range("AF13").Value = 10
$ws.Range("AG13").Value = 35
$ws.Range("AH13").Value = 23
$ws.Range("AI13").Value = 28
$ws.Range("AJ13").Value = 300
# Row 15
$ws.Range("J15").Value = 1.02
$ws.Range("K15").Value = 21
# Row 16
$ws.Range("J16").Value = 1.03
$ws.Range("K16").Value = 17
$ws.Range("L16").Value = 1.14
$ws.Range("M16").Value = 5.5
# Row 17
$ws.Range("G17").Value = 1.44
$ws.Range("I17").Value = 6.25
$ws.Range("J17").Value = 1.03
$ws.Range("K17").Value = 17
$ws.Range("Z17").Value = 17
$ws.Range("AA17").Value = 9
# Row 18
$ws.Range("G18").Value = 2
$ws.Range("I18").Value = 3.7
$ws.Range("R18").Value = 2.1
$ws.Range("S18").Value = 1.67
$ws.Range("U18").Value = 8.5
$ws.Range("AB18").Value = 19
$ws.Range("AF18").Value = 15
# Row 19
$ws.Range("J19").Value = 1.08
$ws.Range("L19").Value = 1.36
# Row 20
$ws.Range("J20").Value = 1.05
$ws.Range("L20").Value = 1.29
# Row 21
$ws.Range("J21").Value = 1.04
$ws.Range("L21").Value = 1.2
# Row 22
$ws.Range("J22").Value = 1.05
$ws.Range("L22").Value = 1.29
# Row 23
$ws.Range("J23").Value = 1.02
$ws.Range("K23").Value = 11
$ws.Range("N23").Value = 1.44
$ws.Range("O23").Value = 2.63
$ws.Range("R23").Value = 1.57
$ws.Range("S23").Value = 2.25
$ws.Range("T23").Value = 11
$ws.Range("U23").Value = 9.5
$ws.Range("V23").Value = 9
$ws.Range("W23").Value = 12
$ws.Range("Y23").Value = 19
$ws.Range("Z23").Value = 21
$ws.Range("AD23").Value = 21
$ws.Range("AE23").Value = 34
$ws.Range("AH23").Value = 34
$ws.Range("AJ23").Value = 126
# Row 25
$ws.Range("I25").Value = 9
$ws.Range("R25").Value = 1.8
$ws.Range("S25").Value = 1.91
$ws.Range("T25").Value = 12
$ws.Range("U25").Value = 8.5
$ws.Range("W25").Value = 8.5
$ws.Range("Z25").Value = 26
$ws.Range("AA25").Value = 15
$ws.Range("AB25").Value = 21
$ws.Range("AD25").Value = 34
$ws.Range("AF25").Value = 29
$ws.Range("AJ25").Value = 600
# Row 27
$ws.Range("H27").Value = 3.35
$ws.Range("I27").Value = 3.95
$ws.Range("K27").Value = 7.1
$ws.Range("L27").Value = 1.31
$ws.Range("M27").Value = 3.15
$ws.Range("N27").Value = 1.91
$ws.Range("O27").Value = 1.8
$ws.Range("P27").Value = 1.44
$ws.Range("Q27").Value = 2.6
$ws.Range("T27").Value = 7
$ws.Range("U27").Value = 8.75
$ws.Range("Z27").Value = 7.1
$ws.Range("AA27").Value = 6.5
$ws.Range("AC27").Value = 70
$ws.Range("AE27").Value = 22
# Row 29
$ws.Range("L29").Value = 1.37
$ws.Range("M29").Value = 2.62
$ws.Range("N29").Value = 2.07
$ws.Range("O29").Value = 1.6
$ws.Range("AB29").Value = 15.5
$ws.Range("AH29").Value = 40
